$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94; existing rows 94-138 shift down to 95-139.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Cells.Item(94, 1).Value = 5
$ws.Cells.Item(94, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(94, 3).Value = "Maule"
$ws.Cells.Item(94, 4).Value = 44510
$ws.Cells.Item(94, 5).Value = 7
$ws.Cells.Item(94, 6).Value = 100112017
$ws.Cells.Item(94, 7).Value = "Apio"
$ws.Cells.Item(94, 8).Value = "Americana (o)"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 500
$ws.Cells.Item(94, 11).Value = 8000
$ws.Cells.Item(94, 12).Value = 8000
$ws.Cells.Item(94, 13).Value = 8000
$ws.Cells.Item(94, 14).Value = "$/docena de matas"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 1333
$ws.Cells.Item(94, 17).Value = 6
$ws.Cells.Item(94, 18).Value = "Hortaliza"
